$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.924.06'
$ws.Range("E2").Value = '  -1.73%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.832.38'
$ws.Range("E3").Value = '  -1.91%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9996'
$ws.Range("E4").Value = '  -0.18%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '240.98'
$ws.Range("E5").Value = '  -0.93%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6870'
$ws.Range("E6").Value = '  -2.73%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9998'
$ws.Range("E7").Value = '  -0.14%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07654'
$ws.Range("E8").Value = '  -2.62%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3048'
$ws.Range("E9").Value = '  -2.63%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '23.57'
$ws.Range("E10").Value = '  -3.77%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07805'
$ws.Range("E11").Value = '  -2.29%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.834.10'
$ws.Range("E12").Value = '  -2.07%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.075'
$ws.Range("E13").Value = '  -2.55%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '90.38'
$ws.Range("E14").Value = '  -3.15%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6770'
$ws.Range("E15").Value = '  -3.18%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.461'
$ws.Range("E16").Value = '  -0.69%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000008274'
$ws.Range("E17").Value = '  -1.31%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '28.926.19'
$ws.Range("E18").Value = '  -1.78%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '242.68'
$ws.Range("E19").Value = '  -3.89%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '2.077.81'
$ws.Range("E20").Value = '  -2.32%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.67'
$ws.Range("E21").Value = '  -3.30%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9998'
$ws.Range("E22").Value = '  -0.13%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.427'
$ws.Range("E23").Value = '  -2.56%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.9997'
$ws.Range("E24").Value = '  -0.18%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1471'
$ws.Range("E25").Value = '  -5.42%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '161.27'
$ws.Range("E26").Value = '  +0.09%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.781'
$ws.Range("E27").Value = '  -2.47%  '
$ws.Range("E28").Value = '  -2.81%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.533'
$ws.Range("E29").Value = '  +2.21%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.213'
$ws.Range("E30").Value = '  -2.56%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.107'
$ws.Range("E31").Value = '  -3.73%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.194'
$ws.Range("E32").Value = '  -0.78%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05117'
$ws.Range("E33").Value = '  -3.61%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7531'
$ws.Range("E34").Value = '  +0.74%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.829'
$ws.Range("E35").Value = '  -3.17%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.144'
$ws.Range("E36").Value = '  -2.29%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.675'
$ws.Range("E37").Value = '  -1.33%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01846'
$ws.Range("E38").Value = '  -1.86%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.227.86'
$ws.Range("E39").Value = '  -3.62%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.695'
$ws.Range("E40").Value = '  -1.87%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9222'
$ws.Range("E41").Value = '  +3.10%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '108.33'
$ws.Range("E42").Value = '  -0.51%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.9992'
$ws.Range("E43").Value = '  -0.16%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.576'
$ws.Range("E44").Value = '  -8.48%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.5170'
$ws.Range("E45").Value = '  -0.15%  '
$ws.Range("B46").Value = 'RocketPoolETH'
$ws.Range("C46").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.977.74'
$ws.Range("E46").Value = '  -2.76%  '
$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.484'
$ws.Range("E47").Value = '  -0.99%  '
$ws.Range("E48").Value = '  -4.80%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '63.94'
$ws.Range("E49").Value = '  -10.41%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.735'
$ws.Range("E50").Value = '  -3.12%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.4188'
$ws.Range("E51").Value = '  -2.73%  '
# End of updates
